$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number format from the last existing data row (A127) down to the
# new rows before filling in values, so the new date cells pick up the same
# custom date style (s="1") instead of creating a duplicate style entry.
$ws.Range("A127").Copy()
$ws.Range("A128:A130").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Add new rows of COVID data for 12/5, 12/6, 12/7 (2020)
# Dates are Excel serial numbers 44170, 44171, 44172
$ws.Range("A128").Value = 44170
$ws.Range("B128").Value = 0
$ws.Range("C128").Value = 0

$ws.Range("A129").Value = 44171
$ws.Range("B129").Value = 0
$ws.Range("C129").Value = 0

$ws.Range("A130").Value = 44172
$ws.Range("B130").Value = 47
$ws.Range("C130").Value = 6

# Update selection to match the final state in the diff
$ws.Range("I137").Select()
